$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 (new row): header line for the Spreadsheet step table ---
$ws.Range("B15").Value = "Spreadsheet SpreadsheetResult calc()"

# --- Row 16 (existing blank row): column headers, default (unstyled) cells ---
$ws.Range("B16").Value = "Step Name"
$ws.Range("C16").Value = "Value"

# --- Row 17 (existing blank row): Step1 ---
$ws.Range("B17").Value = "Step1"
$ws.Range("C20").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = "'= for (int java=0;java<10;java++) {} java.lang.Boolean.TRUE;"

# --- Row 18 (existing blank row): Step2, B gets the bold/wrap style copied from B3 ---
$ws.Range("B3").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "Step2"
$ws.Range("C20").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "'= for (int i=0;i<10;i++) { String java = ""hello"";} java.lang.Boolean.TRUE;"

# --- Row 19: Step3 ---
$ws.Range("B3").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B19").Value = "Step3"
$ws.Range("C20").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C19").Value = "'= if (true) {String java = ""Hello"";} java.lang.Boolean.TRUE;"

# --- Row 20: Step4 (C20 already carried the quote-prefixed style) ---
$ws.Range("B3").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B20").Value = "Step4"
$ws.Range("C20").Value = "'= while (false) {String java = ""Hello"";} java.lang.Boolean.TRUE;"

# --- Row 21: Step5 (C21 already carried the quote-prefixed style) ---
$ws.Range("B3").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").Value = "Step5"
$ws.Range("C21").Value = "'= {String java = ""Hello"";} java.lang.Boolean.TRUE;"

# --- Row 22 (brand new row): Step6 ---
$ws.Range("B22").Value = "Step6"
$ws.Range("C20").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "'while ("

$ws.Application.CutCopyMode = $false

$ws.Range("C22").Select() | Out-Null
